$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.899.18"
$ws.Range("E2").Value = "  +4.48%  "
$ws.Range("D3").Value = "3.361.20"
$ws.Range("E3").Value = "  +5.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'559.72"
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("D6").Value = "'153.40"
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "3.936.47"
$ws.Range("E12").Value = "  +4.83%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'27.23"
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("D16").Value = "62.885.57"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "3.328.00"
$ws.Range("E17").Value = "  +4.13%  "
$ws.Range("D18").Value = "'6.49"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").Value = "'13.84"
$ws.Range("E19").Value = "  +5.73%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "'390.16"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'0.541"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'70.48"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +5.16%  "
$ws.Range("D26").Value = "'8.86"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "0.0₃0975"
$ws.Range("E27").Value = "  +7.36%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'6.59"
$ws.Range("E29").Value = "  +6.40%  "
$ws.Range("E30").Value = "  +4.03%  "
$ws.Range("E31").Value = "  +4.25%  "
$ws.Range("D32").Value = "'23.06"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  +6.77%  "
$ws.Range("D34").Value = "'6.73"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'161.24"
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.49"
$ws.Range("E36").Value = "  +9.42%  "
$ws.Range("E37").Value = "  +11.96%  "
$ws.Range("D38").Value = "'27.10"
$ws.Range("E38").Value = "  +4.87%  "
$ws.Range("D39").Value = "'0.0745"
$ws.Range("E39").Value = "  +4.79%  "
$ws.Range("D40").Value = "2.832.40"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'0.0311"
$ws.Range("E41").Value = "  +8.70%  "
$ws.Range("D42").Value = "'4.32"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "'0.750"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("D44").Value = "'40.78"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").Value = "'22.28"
$ws.Range("E46").Value = "  +8.31%  "
$ws.Range("D47").Value = "3.402.64"
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("D50").Value = "'0.810"
$ws.Range("D51").Value = "'283.17"
$ws.Range("E51").Value = "  +5.31%  "
